# Auto-generated Excel COM-interop script to apply the commit diff
# to Sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR (Brynhildr_Profits workbook).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 525
$ws.Range("J17").Value = 442.22223
$ws.Range("L17").Value = 1326.66669
$ws.Range("N17").Value = -1662.66669
$ws.Range("H112").Value = 1720.4546
$ws.Range("J112").Value = 1635.7894
$ws.Range("L112").Value = 4907.3682
$ws.Range("N112").Value = -7123.3682
$ws.Range("H121").Value = 4284.5625
$ws.Range("J121").Value = 4284.5625
$ws.Range("L121").Value = 12853.6875
$ws.Range("N121").Value = -16347.6875
$ws.Range("H129").Value = 1512.4762
$ws.Range("I129").Value = 1179.8
$ws.Range("K129").Value = 3539.4
$ws.Range("M129").Value = 1460.6
$ws.Range("H132").Value = 9761.846
$ws.Range("I132").Value = 9761.846
$ws.Range("K132").Value = 29285.538
$ws.Range("M132").Value = -26755.538
$ws.Range("H136").Value = 78780
$ws.Range("J136").Value = 78780
$ws.Range("L136").Value = 78780
$ws.Range("N136").Value = -88980
$ws.Range("H137").Value = 6902354.5
$ws.Range("I137").Value = 12501493
$ws.Range("J137").Value = 11107.692
$ws.Range("K137").Value = 37504479
$ws.Range("L137").Value = 33323.076
$ws.Range("M137").Value = -37501929
$ws.Range("N137").Value = -38423.076
$ws.Range("H138").Value = 6386.6807
$ws.Range("I138").Value = 5421.353
$ws.Range("J138").Value = 6933.7
$ws.Range("K138").Value = 16264.059
$ws.Range("L138").Value = 20801.1
$ws.Range("M138").Value = -11124.059
$ws.Range("N138").Value = -31081.1

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 553064.8
$ws.Range("I32").Value = 578056.0600000001
$ws.Range("K32").Value = 578056.0600000001
$ws.Range("M32").Value = -577769.0600000001
$ws.Range("H74").Value = 858548
$ws.Range("I74").Value = 1038222.9
$ws.Range("K74").Value = 1038222.9
$ws.Range("M74").Value = -1037348.9
$ws.Range("H77").Value = 858548
$ws.Range("I77").Value = 1038222.9
$ws.Range("K77").Value = 5191114.5
$ws.Range("M77").Value = -5186746.5
$ws.Range("H132").Value = 6787.8237
$ws.Range("J132").Value = 7646.273
$ws.Range("L132").Value = 22938.819
$ws.Range("N132").Value = -27998.819

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 40000
$ws.Range("I38").Value = 30000
$ws.Range("K38").Value = 30000
$ws.Range("M38").Value = -29584
$ws.Range("H105").Value = 10491.615
$ws.Range("I105").Value = 9126.454
$ws.Range("K105").Value = 9126.454
$ws.Range("M105").Value = -7379.454

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 481391.2
$ws.Range("I31").Value = 614392.1
$ws.Range("J31").Value = 4804.5
$ws.Range("K31").Value = 614392.1
$ws.Range("L31").Value = 4804.5
$ws.Range("M31").Value = -614097.1
$ws.Range("N31").Value = -5394.5
$ws.Range("H34").Value = 481391.2
$ws.Range("I34").Value = 614392.1
$ws.Range("J34").Value = 4804.5
$ws.Range("K34").Value = 614392.1
$ws.Range("L34").Value = 4804.5
$ws.Range("M34").Value = -614190.1
$ws.Range("N34").Value = -5208.5
$ws.Range("H105").Value = 27297
$ws.Range("I105").Value = 26996.666
$ws.Range("K105").Value = 26996.666
$ws.Range("M105").Value = -25249.666
$ws.Range("H134").Value = 4193.684
$ws.Range("I134").Value = 2463.88
$ws.Range("J134").Value = 5545.0938
$ws.Range("K134").Value = 7391.64
$ws.Range("L134").Value = 16635.2814
$ws.Range("M134").Value = -4856.64
$ws.Range("N134").Value = -21705.2814

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1064.8572
$ws.Range("I17").Value = 1000.875
$ws.Range("J17").Value = 1150.1666
$ws.Range("K17").Value = 3002.625
$ws.Range("L17").Value = 3450.4998
$ws.Range("M17").Value = -2833.625
$ws.Range("N17").Value = -3788.4998
$ws.Range("H64").Value = 7249.28
$ws.Range("I64").Value = 3790.8
$ws.Range("K64").Value = 11372.4
$ws.Range("M64").Value = -11102.4
$ws.Range("H67").Value = 7249.28
$ws.Range("I67").Value = 3790.8
$ws.Range("K67").Value = 11372.4
$ws.Range("M67").Value = -10436.4

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2483.0527
$ws.Range("I102").Value = 1828.625
$ws.Range("K102").Value = 1828.625
$ws.Range("M102").Value = -206.625
$ws.Range("H113").Value = 1468.75
$ws.Range("J113").Value = 920
$ws.Range("L113").Value = 920
$ws.Range("N113").Value = -5260
$ws.Range("H126").Value = 2790.9
$ws.Range("I126").Value = 2730.5715
$ws.Range("K126").Value = 8191.7145
$ws.Range("M126").Value = -5721.7145

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6158.9
$ws.Range("I7").Value = 5584.143
$ws.Range("K7").Value = 5584.143
$ws.Range("M7").Value = -5472.143
$ws.Range("H55").Value = 1440.2162
$ws.Range("I55").Value = 1235.3125
$ws.Range("J55").Value = 1596.3334
$ws.Range("K55").Value = 1235.3125
$ws.Range("L55").Value = 1596.3334
$ws.Range("M55").Value = -1062.3125
$ws.Range("N55").Value = -1942.3334
$ws.Range("H93").Value = 2120.35
$ws.Range("I93").Value = 1557.5714
$ws.Range("K93").Value = 1557.5714
$ws.Range("M93").Value = -309.5714
$ws.Range("H123").Value = 100001
$ws.Range("J123").Value = 100001
$ws.Range("L123").Value = 100001
$ws.Range("N123").Value = -109801
$ws.Range("H126").Value = 6158.9
$ws.Range("I126").Value = 5584.143
$ws.Range("K126").Value = 16752.429
$ws.Range("M126").Value = -14282.429
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()  # was -54920, now blank
$ws.Range("H132").Value = 5212480
$ws.Range("I132").Value = 6414206
$ws.Range("K132").Value = 19242618
$ws.Range("M132").Value = -19240088
$ws.Range("H136").Value = 19447242
$ws.Range("I136").Value = 18752664
$ws.Range("K136").Value = 56257992
$ws.Range("M136").Value = -56255442

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2536.2727
$ws.Range("I107").Value = 1157
$ws.Range("K107").Value = 3471
$ws.Range("M107").Value = -1551
$ws.Range("H126").Value = 2615.6
$ws.Range("I126").Value = 2694.5833
$ws.Range("J126").Value = 2299.6667
$ws.Range("K126").Value = 8083.749899999999
$ws.Range("L126").Value = 6899.000100000001
$ws.Range("M126").Value = -5613.749899999999
$ws.Range("N126").Value = -11839.0001
$ws.Range("H132").Value = 5378094
$ws.Range("I132").Value = 8334583.5
$ws.Range("K132").Value = 25003750.5
$ws.Range("M132").Value = -25001220.5
